$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Fix three historical "detect_structure" (Q) values that were miscounted ---
$ws.Range("Q64").Value = 0
$ws.Range("Q70").Value = 0
$ws.Range("Q73").Value = 0

# --- 2) Rows 756 & 757 previously had an empty/placeholder "backup" (R) cell;
#        it now gets an explicit numeric 0 like the rest of the column ---
$ws.Range("R756").Value = 0
$ws.Range("R757").Value = 0

# --- 3) Append 23 new weekly bars (rows 758-780), continuing the series from
#        2024-07-01 through 2024-12-02 ---
$rows = @(
    @(45474, 739.8499755859375, 752, 716.4000244140625, 736.5999755859375, 736.5999755859375, 10947650, 2024, 7, 1, 0, 0, 0, 27, 1, 0, 0),
    @(45481, 742.7999877929688, 746.8499755859375, 695.8499755859375, 713.1500244140625, 713.1500244140625, 8184091, 2024, 7, 8, 0, 0, 0, 28, 0, 0, 0),
    @(45488, 721.3499755859375, 724, 685.5499877929688, 705, 705, 9377360, 2024, 7, 15, 0, 0, 0, 29, 0, 0, 0),
    @(45495, 717, 731.9500122070312, 666, 693.5499877929688, 693.5499877929688, 17709034, 2024, 7, 22, 0, 0, 0, 30, 0, 0, 0),
    @(45502, 698.0499877929688, 748, 680.7000122070312, 722.0499877929688, 722.0499877929688, 17567047, 2024, 7, 29, 0, 0, 0, 31, 0, 0, 0),
    @(45509, 691, 718, 676, 703, 703, 12247517, 2024, 8, 5, 0, 0, 0, 32, 0, 0, 0),
    @(45516, 702, 702, 642.0499877929688, 665.0499877929688, 665.0499877929688, 7596856, 2024, 8, 12, 0, 0, 0, 33, 2, 0, 2),
    @(45523, 665.0999755859375, 725, 664, 703.8499755859375, 703.8499755859375, 9544321, 2024, 8, 19, 0, 0, 0, 34, 0, 0, 0),
    @(45530, 712.5999755859375, 751.3499755859375, 701.9500122070312, 711.1500244140625, 711.1500244140625, 11499939, 2024, 8, 26, 0, 0, 0, 35, 0, 0, 0),
    @(45537, 715.5, 727, 676.25, 709.5499877929688, 709.5499877929688, 12084517, 2024, 9, 2, 0, 0, 0, 36, 0, 0, 0),
    @(45544, 708.25, 775.5999755859375, 705.0999755859375, 768.25, 768.25, 19191156, 2024, 9, 9, 0, 0, 0, 37, 0, 0, 0),
    @(45551, 776, 789.5, 731, 782.3499755859375, 782.3499755859375, 16726832, 2024, 9, 16, 0, 0, 0, 38, 0, 0, 0),
    @(45558, 785.3499755859375, 804.9000244140625, 730, 740.2999877929688, 740.2999877929688, 17687003, 2024, 9, 23, 0, 0, 0, 39, 1, 0, 1),
    @(45565, 735, 736.5999755859375, 691.25, 699.5, 699.5, 9203056, 2024, 9, 30, 0, 0, 0, 40, 0, 0, 0),
    @(45572, 700, 741.4500122070312, 662.2000122070312, 706.7000122070312, 706.7000122070312, 11889887, 2024, 10, 7, 0, 0, 0, 41, 0, 0, 0),
    @(45579, 715.5, 718.4000244140625, 666.0999755859375, 676.8499755859375, 676.8499755859375, 7843730, 2024, 10, 14, 0, 0, 0, 42, 0, 0, 0),
    @(45586, 684, 697.2000122070312, 641.6500244140625, 668.9000244140625, 668.9000244140625, 11442910, 2024, 10, 21, 0, 0, 0, 43, 0, 0, 0),
    @(45593, 672.1500244140625, 694.0999755859375, 635.5499877929688, 682.6500244140625, 682.6500244140625, 9726714, 2024, 10, 28, 0, 0, 0, 44, 2, 0, 0),
    @(45600, 683, 721.8499755859375, 652.0999755859375, 712.0499877929688, 712.0499877929688, 12620526, 2024, 11, 4, 0, 0, 0, 45, 0, 0, 2),
    @(45607, 709.8499755859375, 776.9000244140625, 706.1500244140625, 732, 732, 18551770, 2024, 11, 11, 0, 0, 0, 46, 0, 0, 0),
    @(45614, 736, 746.25, 682.0999755859375, 688.9000244140625, 688.9000244140625, 6865445, 2024, 11, 18, 0, 0, 0, 47, 0, 0, 0),
    @(45621, 705, 711.5, 648.0499877929688, 654.6500244140625, 654.6500244140625, 21040836, 2024, 11, 25, 0, 0, 0, 48, 0, 0, 0),
    @(45628, 655, 680, 640.25, 674.3499755859375, 674.3499755859375, 14334290, 2024, 12, 2, 0, 0, 0, 49, 0, 0, 0)
)

$firstNewRow = 758
$lastNewRow  = $firstNewRow + $rows.Count - 1
$ncols = 17   # columns A..Q

# Column A on the existing rows carries a datetime number format; apply the
# same format to the new rows' A cells before/while writing their values.
$ws.Range("A$($firstNewRow):A$($lastNewRow)").NumberFormat = "YYYY-MM-DD HH:MM:SS"

$data = New-Object 'object[,]' $rows.Count,$ncols
for ($r = 0; $r -lt $rows.Count; $r++) {
    for ($c = 0; $c -lt $ncols; $c++) {
        $data[$r,$c] = $rows[$r][$c]
    }
}
$ws.Range("A$($firstNewRow):Q$($lastNewRow)").Value = $data

# Column R ("backup") is left blank/empty on all the new rows, matching the
# original placeholder state that rows 756/757 used to have.
$ws.Range("R$($firstNewRow):R$($lastNewRow)").Font.Bold = $false
